# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) figures in column F across the
# 展览 / 演出 / 本地生活 sheets and their aggregated mirror in 全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 3347
$wsExpo.Range("F6").Value  = 1131
$wsExpo.Range("F8").Value  = 2168
$wsExpo.Range("F9").Value  = 1138
$wsExpo.Range("F10").Value = 618
$wsExpo.Range("F17").Value = 249
$wsExpo.Range("F22").Value = 277
$wsExpo.Range("F24").Value = 12395
$wsExpo.Range("F25").Value = 12440
$wsExpo.Range("F30").Value = 37
$wsExpo.Range("F31").Value = 411
$wsExpo.Range("F33").Value = 6
$wsExpo.Range("F36").Value = 627

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F9").Value = 43

# --- Sheet "本地生活" (Local Life) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 108

# --- Sheet "全部类型" (All Types, aggregate of the above) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 3347
$wsAll.Range("F7").Value  = 1131
$wsAll.Range("F9").Value  = 2168
$wsAll.Range("F10").Value = 1138
$wsAll.Range("F11").Value = 618
$wsAll.Range("F12").Value = 108
$wsAll.Range("F22").Value = 249
$wsAll.Range("F27").Value = 277
$wsAll.Range("F30").Value = 12395
$wsAll.Range("F31").Value = 12440
$wsAll.Range("F36").Value = 37
$wsAll.Range("F37").Value = 411
$wsAll.Range("F41").Value = 6
$wsAll.Range("F44").Value = 43
$wsAll.Range("F46").Value = 627
